# Update Name of Algo
# Apply updated KNN imputation result values to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = -21.525
$ws.Range("B4").Value = 6.125999999999999
$ws.Range("C4").Value = -12.365
$ws.Range("E4").Value = 12.718

$ws.Range("B5").Value = 6.226

$ws.Range("A7").Value = -21.23

$ws.Range("B8").Value = 6.256

$ws.Range("C9").Value = -11.854

$ws.Range("E12").Value = 13.003

$ws.Range("A16").Value = -21.312
$ws.Range("B16").Value = 5.999000000000001

$ws.Range("C18").Value = -12.732

$ws.Range("E20").Value = 13.198
